$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("profile")

# Insert a new column before column N (the 14th column) so every existing
# column from N onward (pro_soil_taxon, pro_soil_series, ...) shifts one
# place to the right, making room for the new "USDA Soil Order" column.
$ws.Columns("N").Insert()

# Populate the new "pro_usda_soil_order" column: header in row 1, and the
# soil order values for the two existing profile records (rows 4 and 5).
$ws.Cells.Item(1, 14).Value = "pro_usda_soil_order"
$ws.Cells.Item(4, 14).Value = "Alfisols"
$ws.Cells.Item(5, 14).Value = "Ultisols"

# Restore a sensible selection/active-sheet state.
$ws.Range("N6").Select() | Out-Null
$wsMeta = $wb.Worksheets.Item("metadata")
$wsMeta.Activate() | Out-Null
$wsMeta.Range("A4").Select() | Out-Null
